# Update the report title in A1 for KHUNTI - JH23, Jharkhand
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Maker Month Wise Data  of KHUNTI - JH23 , Jharkhand (2022)"

# ---------------------------------------------------------------------------
# Add the two new maker data rows (5 and 6).
# Numeric-looking values must be written as *text* (shared-string) cells, so
# we briefly force a text number format, assign the value, then strip the
# format back off (ClearFormats) so the cell keeps no explicit style - this
# mirrors the workbook's existing cells, which hold numbers as plain text
# without any style index.
# ---------------------------------------------------------------------------

$textRangeA = $ws.Range("A5:A6")
$textRangeB = $ws.Range("C5:O6")
$textRangeA.NumberFormat = "@"
$textRangeB.NumberFormat = "@"

# Row 5 - OKINAWA AUTOTECH PVT LTD
$ws.Range("A5").Value = "1"
$ws.Range("B5").Value = "OKINAWA AUTOTECH PVT LTD"
$ws.Range("C5").Value = "0"
$ws.Range("D5").Value = "0"
$ws.Range("E5").Value = "0"
$ws.Range("F5").Value = "14"
$ws.Range("G5").Value = "7"
$ws.Range("H5").Value = "7"
$ws.Range("I5").Value = "2"
$ws.Range("J5").Value = "6"
$ws.Range("K5").Value = "4"
$ws.Range("L5").Value = "10"
$ws.Range("M5").Value = "4"
$ws.Range("N5").Value = "1"
$ws.Range("O5").Value = "55"

# Row 6 - OLA ELECTRIC TECHNOLOGIES PVT LTD
$ws.Range("A6").Value = "2"
$ws.Range("B6").Value = "OLA ELECTRIC TECHNOLOGIES PVT LTD"
$ws.Range("C6").Value = "0"
$ws.Range("D6").Value = "0"
$ws.Range("E6").Value = "0"
$ws.Range("F6").Value = "0"
$ws.Range("G6").Value = "0"
$ws.Range("H6").Value = "0"
$ws.Range("I6").Value = "0"
$ws.Range("J6").Value = "1"
$ws.Range("K6").Value = "0"
$ws.Range("L6").Value = "0"
$ws.Range("M6").Value = "0"
$ws.Range("N6").Value = "0"
$ws.Range("O6").Value = "1"

# Strip the temporary text format back off so these cells end up with the
# default (no explicit) style, same as the rest of the data rows.
$textRangeA.ClearFormats()
$textRangeB.ClearFormats()

# ---------------------------------------------------------------------------
# Trailing blank separator row, now pushed down to row 11.
# ---------------------------------------------------------------------------
$ws.Range("A11").Font.Bold = $false
$ws.Range("A11").ClearFormats()

# ---------------------------------------------------------------------------
# Column widths - columns A (S No), B (Maker name) and O (Total) need to
# widen to fit the newly added content.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 1.38
$ws.Columns("B").ColumnWidth = 36.62
$ws.Columns("O").ColumnWidth = 2.55
